# Trade #18 (row id 48 / "Trade #" 48 on MarketMaking strategy) closed at
# 2026-02-18 00:12:24 - unknown UNKNOWN +0.000%
#
# Updates the open MarketMaking trade (row 49 on "All Trades", row 20 on
# "MarketMaking") to CLOSED with an early_exit, and rolls the resulting
# P&L through the Summary and Strategy Status roll-up sheets.

$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.68   # Current Capital
$summary.Range("B4").Value = 0.78      # Total P&L $
$summary.Range("B5").Value = 0.34      # Total P&L %
$summary.Range("B6").Value = 46        # Total Trades
$summary.Range("B8").Value = 17        # Losing Trades
$summary.Range("B9").Value = 56.52     # Win Rate %

# --- Strategy Status sheet (MarketMaking row) -----------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.68000000000001  # Capital
$status.Range("D6").Value = 17                  # Trades
$status.Range("E6").Value = -0.13               # P&L $
$status.Range("F6").Value = -0.32               # P&L %
$status.Range("G6").Value = 58.82               # Win Rate %

# --- All Trades sheet (row 49, Trade # 48) --------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G49").Value = 0.26
$allTrades.Range("H49").Value = "CLOSED"
$allTrades.Range("I49").Value = -25.7143
$allTrades.Range("J49").Value = -0.09
$allTrades.Range("K49").Value = 99.68000000000001
$allTrades.Range("L49").Value = "early_exit"
$allTrades.Range("M49").Value = 0.12

# --- MarketMaking strategy sheet (row 20, Trade # 48) ---------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G20").Value = 0.26
$marketMaking.Range("H20").Value = "CLOSED"
$marketMaking.Range("I20").Value = -25.7143
$marketMaking.Range("J20").Value = -0.09
$marketMaking.Range("K20").Value = 99.68000000000001
$marketMaking.Range("P20").Value = "early_exit"
$marketMaking.Range("Q20").Value = 0.12
